# The underlying NATMI TPM computation was redone, which (a) changed most
# of the numeric score columns and (b) dropped "MuSCs" as a possible
# Target cluster (it is still present as a Sending cluster). That turns the
# original 3x3 sending/target cluster cross product (9 data rows) into a
# 3x2 cross product (6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three rows whose "Target cluster" (column D) was MuSCs.
# Deleted bottom-up so the remaining row numbers don't shift underneath us.
$ws.Rows("10:10").Delete() | Out-Null
$ws.Rows("7:7").Delete() | Out-Null
$ws.Rows("4:4").Delete() | Out-Null

# The six surviving rows (now rows 2-7) keep the same column layout (A-T)
# but get refreshed Sending/Target cluster labels plus the recomputed
# TPM-derived metrics.

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfrsf14"
$ws.Range("C2").Value = "Btla"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.558321333333333
$ws.Range("H2").Value = 10.674964
$ws.Range("I2").Value = 0.3039644761000113
$ws.Range("J2").Value = 0.3039644761000113
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3107469999999999
$ws.Range("N2").Value = 0.9322409999999999
$ws.Range("O2").Value = 0.9278633407583023
$ws.Range("P2").Value = 0.9278633407583025
$ws.Range("Q2").Value = 1.105737679369333
$ws.Range("R2").Value = 9.951639114323998
$ws.Range("S2").Value = 0.2820374942660037
$ws.Range("T2").Value = 0.2820374942660036

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfrsf14"
$ws.Range("C3").Value = "Btla"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.558321333333333
$ws.Range("H3").Value = 10.674964
$ws.Range("I3").Value = 0.3039644761000113
$ws.Range("J3").Value = 0.3039644761000113
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.024159
$ws.Range("N3").Value = 0.072477
$ws.Range("O3").Value = 0.07213665924169768
$ws.Range("P3").Value = 0.07213665924169768
$ws.Range("Q3").Value = 0.085965485092
$ws.Range("R3").Value = 0.773689365828
$ws.Range("S3").Value = 0.02192698183400768
$ws.Range("T3").Value = 0.02192698183400767

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfrsf14"
$ws.Range("C4").Value = "Btla"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.383140666666667
$ws.Range("H4").Value = 16.149422
$ws.Range("I4").Value = 0.4598470400038817
$ws.Range("J4").Value = 0.4598470400038817
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3107469999999999
$ws.Range("N4").Value = 0.9322409999999999
$ws.Range("O4").Value = 0.9278633407583023
$ws.Range("P4").Value = 0.9278633407583025
$ws.Range("Q4").Value = 1.672794812744666
$ws.Range("R4").Value = 15.055153314702
$ws.Range("S4").Value = 0.4266752107758184
$ws.Range("T4").Value = 0.4266752107758184

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfrsf14"
$ws.Range("C5").Value = "Btla"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.383140666666667
$ws.Range("H5").Value = 16.149422
$ws.Range("I5").Value = 0.4598470400038817
$ws.Range("J5").Value = 0.4598470400038817
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.024159
$ws.Range("N5").Value = 0.072477
$ws.Range("O5").Value = 0.07213665924169768
$ws.Range("P5").Value = 0.07213665924169768
$ws.Range("Q5").Value = 0.130051295366
$ws.Range("R5").Value = 1.170461658294
$ws.Range("S5").Value = 0.03317182922806334
$ws.Range("T5").Value = 0.03317182922806333

# Row 6: MuSCs -> ECs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Tnfrsf14"
$ws.Range("C6").Value = "Btla"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.764910333333333
$ws.Range("H6").Value = 8.294730999999999
$ws.Range("I6").Value = 0.2361884838961071
$ws.Range("J6").Value = 0.236188483896107
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3107469999999999
$ws.Range("N6").Value = 0.9322409999999999
$ws.Range("O6").Value = 0.9278633407583023
$ws.Range("P6").Value = 0.9278633407583025
$ws.Range("Q6").Value = 0.859187591352333
$ws.Range("R6").Value = 7.732688322170998
$ws.Range("S6").Value = 0.2191506357164804
$ws.Range("T6").Value = 0.2191506357164804

# Row 7: MuSCs -> FAPs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Tnfrsf14"
$ws.Range("C7").Value = "Btla"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.764910333333333
$ws.Range("H7").Value = 8.294730999999999
$ws.Range("I7").Value = 0.2361884838961071
$ws.Range("J7").Value = 0.236188483896107
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.024159
$ws.Range("N7").Value = 0.072477
$ws.Range("O7").Value = 0.07213665924169768
$ws.Range("P7").Value = 0.07213665924169768
$ws.Range("Q7").Value = 0.06679746874299999
$ws.Range("R7").Value = 0.6011772186869999
$ws.Range("S7").Value = 0.01703784817962668
$ws.Range("T7").Value = 0.01703784817962667
